# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 306 of the single
# worksheet ("Fruta, Feria Lagunitas de Puerto Montt - Naranja"). This
# pushes the former rows 306-337 down to 307-338 (dimension grows from
# A1:T337 to A1:T338) and the new row is populated with the latest
# Valencia/Primera price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 306 - shifts rows 306:337 down to 307:338 and
# carries the date-format style from the cell above into the new D306.
$ws.Rows(306).Insert()

$ws.Range("A306").Value = 4
$ws.Range("B306").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C306").Value = "Los Lagos"
$ws.Range("D306").Value2 = 44578
$ws.Range("E306").Value = 10
$ws.Range("F306").Value = "Fruta"
$ws.Range("G306").Value = 100102
$ws.Range("H306").Value = "Cítricos"
$ws.Range("I306").Value = 100102005
$ws.Range("J306").Value = "Naranja"
$ws.Range("K306").Value = "Valencia"
$ws.Range("L306").Value = "Primera"
$ws.Range("M306").Value = 400
$ws.Range("N306").Value = 17000
$ws.Range("O306").Value = 18000
$ws.Range("P306").Value = 17500
$ws.Range("Q306").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R306").Value = "Región de O'Higgins"
$ws.Range("S306").Value = 1167
$ws.Range("T306").Value = 15
